$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.664.68"
$ws.Range("E2").Value = "  -1.86%  "

$ws.Range("D3").Value = "1.797.83"
$ws.Range("E3").Value = "  -1.54%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5879"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.06%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06799"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07529"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("D12").Value = "1.804.66"
$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.787"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6196"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").Value = "2.042.05"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009110"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.11%  "

$ws.Range("D18").Value = "28.632.62"
$ws.Range("E18").Value = "  -1.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.485"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "210.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.828"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.14%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.964"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1268"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.425"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06129"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.425"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.823"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.784"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.738"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.052"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6439"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.498"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.714"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.518"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01698"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("D41").Value = "1.144.97"
$ws.Range("E41").Value = "  -6.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8863"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.05%  "

$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("D45").Value = "1.949.28"
$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.19%  "

$ws.Range("E47").Value = "  -4.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.597"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.333"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.14%  "

$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4481"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.52%  "
